$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 44544
$ws.Range('M2').Value = 400
$ws.Range('O2').Value = 5500
$ws.Range('P2').Value = 5250
$ws.Range('R2').Value = 'Región del Maule'
$ws.Range('S2').Value = 3500
$ws.Range('D3').Value = 44537
$ws.Range('D4').Value = 44169
$ws.Range('N4').Value = 5500
$ws.Range('O4').Value = 6000
$ws.Range('P4').Value = 5750
$ws.Range('R4').Value = 'Provincia de Curicó'
$ws.Range('S4').Value = 3833
$ws.Range('D5').Value = 44553
$ws.Range('D6').Value = 44551
$ws.Range('D7').Value = 44547
$ws.Range('L7').Value = 'Primera'
$ws.Range('M7').Value = 400
$ws.Range('N7').Value = 5000
$ws.Range('O7').Value = 5500
$ws.Range('P7').Value = 5250
$ws.Range('R7').Value = 'Región del Maule'
$ws.Range('S7').Value = 3500
$ws.Range('D8').Value = 44530
$ws.Range('M8').Value = 160
$ws.Range('N8').Value = 3600
$ws.Range('O8').Value = 3700
$ws.Range('P8').Value = 3650
$ws.Range('S8').Value = 3650
$ws.Range('D9').Value = 44159
$ws.Range('L9').Value = 'Segunda'
$ws.Range('N9').Value = 6500
$ws.Range('O9').Value = 7000
$ws.Range('P9').Value = 6750
$ws.Range('S9').Value = 4500
$ws.Range('D10').Value = 44516
$ws.Range('M10').Value = 80
$ws.Range('N10').Value = 3700
$ws.Range('O10').Value = 3800
$ws.Range('P10').Value = 3750
$ws.Range('S10').Value = 3750
$ws.Range('D11').Value = 44166
$ws.Range('M11').Value = 200
$ws.Range('N11').Value = 6000
$ws.Range('O11').Value = 6500
$ws.Range('P11').Value = 6250
$ws.Range('Q11').Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range('R11').Value = 'Provincia de Curicó'
$ws.Range('S11').Value = 4167
$ws.Range('T11').Value = 1.5
$ws.Range('D12').Value = 44519
$ws.Range('M12').Value = 200
$ws.Range('D13').Value = 44533
$ws.Range('N13').Value = 3500
$ws.Range('O13').Value = 3600
$ws.Range('P13').Value = 3550
$ws.Range('Q13').Value = '$/kilo'
$ws.Range('S13').Value = 3550
$ws.Range('T13').Value = 1
$ws.Range('D14').Value = 44523
$ws.Range('M14').Value = 300
$ws.Range('N14').Value = 3700
$ws.Range('O14').Value = 3800
$ws.Range('P14').Value = 3750
$ws.Range('Q14').Value = '$/kilo'
$ws.Range('R14').Value = 'Región del Maule'
$ws.Range('S14').Value = 3750
$ws.Range('T14').Value = 1
$ws.Range('D15').Value = 44176
$ws.Range('M15').Value = 300
$ws.Range('N15').Value = 5000
$ws.Range('O15').Value = 6000
$ws.Range('P15').Value = 5500
$ws.Range('S15').Value = 3667
$ws.Range('L16').Value = 'Primera'
$ws.Range('N16').Value = 7000
$ws.Range('O16').Value = 7000
$ws.Range('P16').Value = 7000
$ws.Range('S16').Value = 4667
$ws.Range('D17').Value = 44162
$ws.Range('L17').Value = 'Segunda'
$ws.Range('M17').Value = 100
$ws.Range('N17').Value = 6500
$ws.Range('O17').Value = 6500
$ws.Range('P17').Value = 6500
$ws.Range('Q17').Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range('R17').Value = 'Provincia de Curicó'
$ws.Range('S17').Value = 4333
$ws.Range('T17').Value = 1.5
